$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes existing rows 3:30 down to 4:31,
# extending the used range to A1:R31) and fill it with this week's new
# data point for the weekly price series.
$ws.Rows("3:3").Insert()

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = "2022-02-24"
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = 100112052
$ws.Cells.Item(3, 7).Value = "Albahaca"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 200
$ws.Cells.Item(3, 11).Value = 2500
$ws.Cells.Item(3, 12).Value = 3000
$ws.Cells.Item(3, 13).Value = 2750
$ws.Cells.Item(3, 14).Value = "`$/paquete"
$ws.Cells.Item(3, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(3, 16).Value = 2750
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"
